$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume data to the latest scrape values.
# Force text storage (matches the original inline-string cell type)
# so numeric-looking strings like "1.00" or "315.60" keep their
# exact formatting instead of being coerced to doubles.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.194.14"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.60%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.598.79"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +3.38%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.60"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.42%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.86"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +4.02%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.578"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.45%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.12%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.539"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.07%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.90"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.01%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.40%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.54"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.36%  "

$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.108"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.46%  "

$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.808.21"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -3.28%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.596.77"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.33%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.29"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.71%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.850"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.32%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "41.973.70"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.58%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.87"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.82"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.77%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0969"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.19%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.69"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.59%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "254.03"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.39%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.35%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.10"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +4.27%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "27.35"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.43%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.12%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.42"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.34%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "41.22"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.40%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.34"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.08%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.76%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "156.02"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.52%  "

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +6.02%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.80%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +4.33%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +3.30%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.94"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.10%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.19%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +9.66%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.00%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "23.08"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.24%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.99"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +5.92%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.27%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.36%  "

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.61%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.016.54"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.28%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.79%  "

$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "BitcoinSV"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "83.23"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.55%  "

$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.197"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +5.13%  "

$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "ordi"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "74.79"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.55%  "

$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "104.00"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.07%  "
